$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Ensure the price/volume columns stay text-typed so numeric-looking
# strings (e.g. '1.000', '0.000009912') are not reinterpreted as numbers.
$ws.Range("D2:E51").NumberFormat = "@"

$ws.Range("D2").Value = '29.023.23'
$ws.Range("E2").Value = '  +0.15%  '
$ws.Range("D3").Value = '1.916.66'
$ws.Range("E3").Value = '  +0.96%  '
$ws.Range("E4").Value = '  -0.44%  '
$ws.Range("D5").Value = '324.97'
$ws.Range("E5").Value = '  -0.59%  '
$ws.Range("D6").Value = '0.9999'
$ws.Range("E6").Value = '  -0.39%  '
$ws.Range("D7").Value = '0.4603'
$ws.Range("E7").Value = '  -0.17%  '
$ws.Range("D8").Value = '0.3872'
$ws.Range("E8").Value = '  -1.37%  '
$ws.Range("D9").Value = '0.07831'
$ws.Range("E9").Value = '  -0.66%  '
$ws.Range("D10").Value = '0.9911'
$ws.Range("E10").Value = '  -0.17%  '
$ws.Range("D11").Value = '21.95'
$ws.Range("E11").Value = '  -1.00%  '
$ws.Range("D12").Value = '1.885.40'
$ws.Range("E12").Value = '  -2.90%  '
$ws.Range("D13").Value = '5.781'
$ws.Range("E13").Value = '  +0.94%  '
$ws.Range("D14").Value = '7.021'
$ws.Range("E14").Value = '  -0.66%  '
$ws.Range("D15").Value = '0.07055'
$ws.Range("E15").Value = '  +1.52%  '
$ws.Range("D16").Value = '87.35'
$ws.Range("E16").Value = '  -1.14%  '
$ws.Range("E17").Value = '  -0.38%  '
$ws.Range("D18").Value = '0.000009912'
$ws.Range("E18").Value = '  -1.01%  '
$ws.Range("D19").Value = '17.05'
$ws.Range("E19").Value = '  +0.33%  '
$ws.Range("D20").Value = '1.000'
$ws.Range("E20").Value = '  -0.35%  '
$ws.Range("D21").Value = '29.046.68'
$ws.Range("E21").Value = '  +0.11%  '
$ws.Range("D22").Value = '5.384'
$ws.Range("E22").Value = '  +1.07%  '
$ws.Range("D23").Value = '11.12'
$ws.Range("E23").Value = '  +0.98%  '
$ws.Range("D24").Value = '2.135.40'
$ws.Range("E24").Value = '  -1.18%  '
$ws.Range("D25").Value = '2.082'
$ws.Range("E25").Value = '  +0.72%  '
$ws.Range("D26").Value = '156.09'
$ws.Range("E26").Value = '  -0.16%  '
$ws.Range("D27").Value = '19.32'
$ws.Range("E27").Value = '  -0.11%  '
$ws.Range("D28").Value = '5.847'
$ws.Range("E28").Value = '  -2.45%  '
$ws.Range("D29").Value = '118.06'
$ws.Range("E29").Value = '  +0.06%  '
$ws.Range("D30").Value = '1.863'
$ws.Range("E30").Value = '  -3.55%  '
$ws.Range("D31").Value = '0.09309'
$ws.Range("E31").Value = '  -0.43%  '
$ws.Range("D32").Value = '0.8795'
$ws.Range("E32").Value = '  -3.79%  '
$ws.Range("D33").Value = '5.193'
$ws.Range("E33").Value = '  -2.53%  '
$ws.Range("D34").Value = '1.311'
$ws.Range("E34").Value = '  -2.41%  '
$ws.Range("D35").Value = '3.127'
$ws.Range("E35").Value = '  -4.85%  '
$ws.Range("E36").Value = '  -0.65%  '
$ws.Range("D38").Value = '0.02087'
$ws.Range("E38").Value = '  -0.03%  '
$ws.Range("D39").Value = '0.9988'
$ws.Range("D40").Value = '7.649'
$ws.Range("E40").Value = '  -2.07%  '
$ws.Range("D41").Value = '0.5678'
$ws.Range("E41").Value = '  -0.44%  '
$ws.Range("D42").Value = '0.1808'
$ws.Range("E42").Value = '  +1.39%  '
$ws.Range("D43").Value = '0.000002995'
$ws.Range("E43").Value = '  +86.43%  '
$ws.Range("D44").Value = '9.678'
$ws.Range("E44").Value = '  -1.75%  '
$ws.Range("D45").Value = '11.79'
$ws.Range("E45").Value = '  -1.82%  '
$ws.Range("D46").Value = '2.196'
$ws.Range("E46").Value = '  -3.20%  '
$ws.Range("D47").Value = '0.5315'
$ws.Range("E47").Value = '  -1.05%  '
$ws.Range("E48").Value = '  -1.71%  '
$ws.Range("D49").Value = '2.565'
$ws.Range("E49").Value = '  +1.24%  '
$ws.Range("D50").Value = '1.835'
$ws.Range("E50").Value = '  -0.90%  '
$ws.Range("D51").Value = '112.53'
$ws.Range("E51").Value = '  -0.29%  '
